$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.45
$ws.Range("I2").Value = 6.5
$ws.Range("L2").Value = 6
$ws.Range("X2").Value = 7.5
$ws.Range("AB2").Value = 23
$ws.Range("AE2").Value = 17
$ws.Range("AJ2").Value = 19
$ws.Range("AL2").Value = 41
